$d = $word.ActiveDocument

# XML-package wrapper helper text (pkg:package) used with Range.InsertXML so we can
# specify exact run boundaries (InsertXML does not auto-coalesce adjacent runs that
# plain text edits would merge together).
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$langRpr = '<w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr>'

# --- 1. Remove the whole "License Information" Heading2 paragraph -----------------
$d.Content.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$licPara = $d.Paragraphs(4)
$licPara.Range.Delete()

# --- 2. Paragraph that used to read "बाइबल कोश (टिंडेल) (Hindi) is based on: ..." --
# becomes just the leading blank run + a bold "Aquifer Open Bible Dictionary" run.
$p5 = $d.Paragraphs(4)
$p5xml = $pkgHeader + '<w:p>' + `
    '<w:r>' + $langRpr + '</w:r>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr><w:t>Aquifer Open Bible Dictionary</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$d.Range($p5.Range.Start, $p5.Range.End - 1).InsertXML($p5xml)

# --- 3. Paragraph "This PDF version is provided under the same license." becomes the
# new adaptation/licensing statement split across several runs.
$p6 = $d.Paragraphs(5)
$p6xml = $pkgHeader + '<w:p>' + `
    '<w:r>' + $langRpr + '<w:t xml:space="preserve">This work is an adaptation of </w:t></w:r>' + `
    '<w:r>' + $langRpr + '<w:t>Tyndale Open Bible Dictionary</w:t></w:r>' + `
    '<w:r>' + $langRpr + '<w:t xml:space="preserve"> &#169; 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, </w:t></w:r>' + `
    '<w:r>' + $langRpr + '<w:t>Aquifer Open Bible Dictionary</w:t></w:r>' + `
    '<w:r>' + $langRpr + '<w:t>, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$d.Range($p6.Range.Start, $p6.Range.End - 1).InsertXML($p6xml)

# --- 4. Brand-new paragraph inserted right after it, describing the translations. ---
$p6 = $d.Paragraphs(5)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs(6)
$p7xml = $pkgHeader + '<w:p>' + `
    '<w:r>' + $langRpr + '</w:r>' + `
    '<w:r>' + $langRpr + '<w:t>This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (&#1593;&#1585;&#1576;&#1610;), French (Fran&#231;ais), Hindi (&#2361;&#2367;&#2306;&#2342;&#2368;), Indonesian (Bahasa Indonesia), Portuguese (Portugu&#234;s), Russian (&#1056;&#1091;&#1089;&#1089;&#1082;&#1080;&#1081;), Spanish (Espa&#241;ol), Swahili (Kiswahili), and Simplified Chinese (&#31616;&#20307;&#20013;&#25991;).</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$d.Range($p7.Range.Start, $p7.Range.End).InsertXML($p7xml)

Write-Output "done"
